$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 ("time_taken"), reusing the same formatting
# (bold, bordered, centered) as the other header cells by copying E1's
# format onto F1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Add the time_taken values for each data row (no special style, matching
# the plain formatting used by the other data cells in column E).
$ws.Range("F2").Value = "2021-10-05 10:51:10.219809"
$ws.Range("F3").Value = "2021-10-05 10:51:10.219819"
$ws.Range("F4").Value = "2021-10-05 10:51:10.219823"
